{"js": "// The document's title, author, and abstract paragraphs each have their\n// text split across many single-word/single-space runs, e.g.\n//   \"Questions:\" + \" \" + \"Introduction\" + \" \" + \"to\" + \" \" + \"sigma\" + \" \" + \"notation\"\n// Consolidate each of those paragraphs down to a single run holding the\n// same, unchanged sentence as one string \u2014 matching the target OOXML.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// Map style name -> full replacement text (identical wording, just no\n// longer split across many runs).\nconst targetsByStyle = {\n  \"Title\": \"Questions: Introduction to sigma notation\",\n  \"Author\": \"Ifan Howells-Baines, Mark Toner\",\n  \"Abstract\": \"Questions relating to the guide on introduction to sigma notation.\"\n};\n\nfor (const paragraph of paragraphs.items) {\n  const newText = targetsByStyle[paragraph.style];\n  if (newText !== undefined) {\n    paragraph.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's title, author, and abstract paragraphs each have their\n# text split across many single-word/single-space runs, e.g.\n#   \"Questions:\" + \" \" + \"Introduction\" + \" \" + \"to\" + \" \" + \"sigma\" + \" \" + \"notation\"\n# Consolidate each of those paragraphs down to a single run holding the\n# same, unchanged sentence as one string -- matching the target OOXML.\n#\n# Paragraphs are located by style (Title / Author / Abstract) rather than\n# a hard-coded index, and each Find/Replace is scoped to that paragraph's\n# own Range so the unrelated, later occurrence of \"Ifan Howells-Baines,\n# Mark Toner\" inside the \"Version history and licensing\" paragraph is\n# left untouched.\n\n$d = $word.ActiveDocument\n\n$targetsByStyle = @{\n    \"Title\"    = \"Questions: Introduction to sigma notation\"\n    \"Author\"   = \"Ifan Howells-Baines, Mark Toner\"\n    \"Abstract\" = \"Questions relating to the guide on introduction to sigma notation.\"\n}\n\nforeach ($para in $d.Paragraphs) {\n    $styleName = $para.Range.Style.NameLocal\n    if ($targetsByStyle.ContainsKey($styleName)) {\n        $newText = $targetsByStyle[$styleName]\n        $range = $para.Range\n        # Exclude the trailing paragraph mark from the replace range.\n        $scope = $d.Range($range.Start, $range.End - 1)\n        $scope.Find.Execute($scope.Text, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    }\n}\n"}
